# Update the "Estado de Cuenta" (account statement) detail table: the previous
# set of mora periods (2312 .. 2503) is replaced by a refreshed set where the
# periods are listed in reverse order (most recent period first: 2503 down to
# 2312), the base salary (column G) is updated from 1,160,000 to 1,800,000 for
# every detail row, and the partial-month "Valor Mora" (27840) moves from the
# first detail row to the last one (column F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31)
$periods = @("2503","2502","2501","2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401","2312")
$moraValues = @(27840,46400,46400,46400,46400,46400,46400,46400,46400,46400,46400,46400,46400,46400,46400,46400)
$baseSalary = 1800000

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("F$r").Value = $moraValues[$i]
    $ws.Range("G$r").Value = $baseSalary
}
